# refactor: remove supplier references from barang management
#
# The "supplier_id" column (column B) is no longer used by barang
# management, so drop it entirely. Excel shifts barang_kode,
# barang_nama, harga_beli and harga_jual one column to the left
# (B:E) and the shared-strings table loses the now-unused
# "supplier_id" entry automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the supplier_id column (B) — cells to the right shift left.
$ws.Columns("B").Delete() | Out-Null

# Match the author's final cursor position in the saved file.
$ws.Range("C11").Select() | Out-Null
